{"js": "// Update the \"Artificial Intelligence\" policy paragraph:\n//  - rewrite the paragraph text to the new policy wording\n//  - italicize the standalone word \"not\" (\"...Artificial intelligence is *not* permitted...\")\n\nconst body = context.document.body;\n\n// Locate the paragraph by a distinctive substring of its ORIGINAL text so the\n// script doesn't depend on a hard-coded paragraph index.\nconst anchorText = \"Based on the assignments in this course and our specified learning outcomes\";\nconst hits = body.search(anchorText, { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not locate the Artificial Intelligence policy paragraph.\");\n}\n\nconst targetParagraph = hits.items[0].paragraphs.getFirst();\ntargetParagraph.load(\"text\");\nawait context.sync();\n\n// The new paragraph text (same content the commit replaces the old run with).\nconst seg1 = \"Based on the specific learning outcomes and assignments in this course, artificial intelligence is permitted on the following: all weekly assignments. Artificial intelligence is\";\nconst seg2 = \" \";\nconst seg3 = \"not\";\nconst seg4 = \" \";\nconst seg5 = \"permitted in tests and practicums. See each assignment, quiz, or exam instructions for more information about what artificial intelligence tools are permitted and to what extent, as well as citation requirements. If no instructions are provided for a specific assignment, then no use of any artificial intelligence tool is permitted. Any AI use beyond that which is detailed in course assignments is explicitly prohibited except when documented permission is granted.\";\n\nconst newText = seg1 + seg2 + seg3 + seg4 + seg5;\n\n// Replace the whole paragraph's text in one shot (clears any pre-existing\n// per-run formatting so the base text comes back as plain runs).\ntargetParagraph.insertText(newText, \"Replace\");\nawait context.sync();\n\n// Re-search, now scoped to this paragraph, for the standalone word \"not\" and\n// italicize it. matchWholeWord avoids hitting \"not\" inside some other word.\nconst notHits = targetParagraph.search(seg3, { matchCase: true, matchWholeWord: true });\nnotHits.load(\"items\");\nawait context.sync();\n\nif (notHits.items.length === 0) {\n  throw new Error('Could not find the word \"not\" to italicize.');\n}\n\n// There should be exactly one \"not\" in the rewritten paragraph; italicize it.\nnotHits.items[0].font.set({ italic: true });\nawait context.sync();\n", "ps1": "# Update the \"Artificial Intelligence\" policy paragraph:\n#  - rewrite the paragraph text to the new policy wording\n#  - italicize the standalone word \"not\" (\"...Artificial intelligence is *not* permitted...\")\n\n$d = $word.ActiveDocument\n\n# Locate the paragraph by a distinctive prefix of its ORIGINAL text so the\n# script doesn't depend on a hard-coded paragraph index.\n$target = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs.Item($i)\n  if ($p.Range.Text -like \"Based on the assignments in this course and our specified learning outcomes*\") {\n    $target = $p\n    break\n  }\n}\n\nif ($target -eq $null) {\n  throw \"Could not locate the Artificial Intelligence policy paragraph.\"\n}\n\n# The new paragraph text (same content the commit replaces the old run with).\n$seg1 = \"Based on the specific learning outcomes and assignments in this course, artificial intelligence is permitted on the following: all weekly assignments. Artificial intelligence is\"\n$seg2 = \" \"\n$seg3 = \"not\"\n$seg4 = \" \"\n$seg5 = \"permitted in tests and practicums. See each assignment, quiz, or exam instructions for more information about what artificial intelligence tools are permitted and to what extent, as well as citation requirements. If no instructions are provided for a specific assignment, then no use of any artificial intelligence tool is permitted. Any AI use beyond that which is detailed in course assignments is explicitly prohibited except when documented permission is granted.\"\n\n$newText = $seg1 + $seg2 + $seg3 + $seg4 + $seg5\n\n# Replace the whole paragraph's text in one shot (the trailing paragraph mark\n# is implicit in Range.Text and must not be included here).\n$target.Range.Text = $newText\n\n# Re-fetch the paragraph (its Range grew) and find the standalone word \"not\"\n# within it, then italicize just that word.\n$target2 = $d.Paragraphs.Item($i)\n$searchRange = $target2.Range\n$searchRange.Find.ClearFormatting()\n$searchRange.Find.Text = $seg3\n$searchRange.Find.MatchCase = $true\n$searchRange.Find.MatchWholeWord = $true\n$found = $searchRange.Find.Execute()\n\nif (-not $found) {\n  throw 'Could not find the word \"not\" to italicize.'\n}\n\n$searchRange.Font.Italic = 1\n"}
